# edit.ps1
#
# Applies the "new quest images, english update" diff to the rules
# paragraphs of kta/pravidla.docx:
#
#   1) "Kradeni" bullet - join the sentence that was split around a
#      stray _GoBack bookmark back into continuous prose (bookmark
#      removed).
#   2) "Jizda" bullet - rewrite the car-ramming dice rule, splitting it
#      into two runs.
#   3) the old "po poskozeni jineho auta ..." sub-bullet is replaced by
#      three new paragraphs describing the new damage/police rules.
#   4) "Strelba" - the old "Otoc kartu ..." bullet is replaced with the
#      "Je-li v aute hrac ..." wording (now with an extra clause).
#   5) the (now-orphaned) original "Je-li v aute hrac ..." bullet is
#      replaced with "Pokud poskodis auto, ziskavas policajta", and the
#      _GoBack bookmark moves here.
#   6) a <w:lastRenderedPageBreak/> is added right before the manual
#      page break run near the end of the document.
#
# Implementation notes:
#  - Plain wording tweaks that stay inside a contiguous run/sentence use
#    Range.Find.Execute (wdReplaceOne), which is also how the stray
#    bookmark from step 1 gets merged away/dropped automatically.
#  - Structural edits (new paragraphs, moved bookmarks, runs with
#    different splits) use Range.InsertXML with a minimal flat-OPC
#    WordprocessingML payload, which replaces the exact contents of the
#    addressed Range. A *fresh* Range object (via $d.Range(start,end))
#    is used for InsertXML, since re-using the same Find-mutated Range
#    instance does not reliably replace in place.

$d = $word.ActiveDocument

function Get-MatchRange([string]$searchText) {
    # Locates searchText (case sensitive) anywhere in the document and
    # returns a fresh Range covering exactly that text.
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return $d.Range($r.Start, $r.End)
}

function Get-MatchParagraphRange([string]$searchText) {
    # Locates searchText (case sensitive) and returns a fresh Range
    # covering the whole enclosing paragraph (including its end-of
    # paragraph mark), so InsertXML can replace pPr/runs wholesale.
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $r.Expand(4) | Out-Null   # wdParagraph
    return $d.Range($r.Start, $r.End)
}

function New-FlatOpcBody([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) "Kradeni": merge "umist" + "i na nej svou figurku" (previously
#    split by a <w:bookmarkStart/End w:name="_GoBack"/> pair) into one
#    continuous sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "umíst" + "i na něj svou figurku",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "umísti na něj svou figurku",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) "Jizda": "... hoď kostkou" -> "..., " / "hoď tolika kostkami, o
#    kolik je tvá rychlost větší než u cílového auta" (two runs).
# ---------------------------------------------------------------------
$rng = Get-MatchRange "poté můžeš nabourat jiné auto, za každou svou aktuální rychlost hoď kostkou"
$rng.InsertXML((New-FlatOpcBody (
    '<w:p><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">poté můžeš nabourat jiné auto, </w:t></w:r>' +
    '<w:r><w:t>hoď tolika kostkami, o kolik je tvá rychlost větší než u cílového auta</w:t></w:r>' +
    '</w:p>'
)))

# ---------------------------------------------------------------------
# 3) Replace the "po poškození jiného auta otoč kartu ... (max. 5)"
#    sub-bullet paragraph with three new paragraphs.
# ---------------------------------------------------------------------
$rng = Get-MatchParagraphRange "po poškození jiného auta otoč kartu a urči policajty, kteří tě začnou stíhat"
$rng.InsertXML((New-FlatOpcBody (
    '<w:p><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:spacing w:after="0"/><w:ind w:left="708"/></w:pPr>' +
    '<w:r><w:t>za každé udělené poškození můžeš zaplatit jednu rychlost a přidat jedno poškození na sebe i na nabourané auto</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:spacing w:after="0"/><w:ind w:left="708"/></w:pPr>' +
    '<w:r><w:t>pokud udělíš aspoň jedno poškození, získáváš policajta</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="992"/></w:pPr></w:p>'
)))

# ---------------------------------------------------------------------
# 5) Replace the (still unique at this point) "Je-li v autě hráč, může
#    střelbu opětovat ..." bullet with "Pokud poškodíš auto, získáváš
#    policajta", moving the _GoBack bookmark here. Done *before* step 4
#    below, which would otherwise re-introduce a duplicate "Je-li v
#    autě hráč ..." sentence and make this search ambiguous.
# ---------------------------------------------------------------------
$rng = Get-MatchParagraphRange "Je-li v autě hráč, může střelbu opětovat (1 žeton, 1 poškození), bez policajtů"
$rng.InsertXML((New-FlatOpcBody (
    '<w:p><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t>Pokud poškodíš auto, získáváš policajta</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:br/></w:r></w:p>'
)))

# ---------------------------------------------------------------------
# 4) Replace "Otoč kartu a urči policajty, kteří tě začnou stíhat (max.
#    5), při střelbě přímo na hráče automaticky 1" with "Je-li v autě
#    hráč, může střelbu opětovat (1 žeton, 1 poškození), bez policajtů"
#    + " anebo za žeton rychlosti poškození zrušit".
# ---------------------------------------------------------------------
$rng = Get-MatchParagraphRange "Otoč kartu a urči policajty, kteří tě začnou stíhat"
$rng.InsertXML((New-FlatOpcBody (
    '<w:p><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t>Je-li v autě hráč, může střelbu opětovat (1 žeton, 1 poškození), bez policajtů</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> anebo za žeton rychlosti poškození zrušit</w:t></w:r></w:p>'
)))

# ---------------------------------------------------------------------
# 6) Add <w:lastRenderedPageBreak/> right before the manual page-break
#    run near the end of the document.
# ---------------------------------------------------------------------
$pageBreakPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.WordOpenXML -like "*<w:br w:type=*") {
        $pageBreakPara = $p
        break
    }
}
if ($pageBreakPara -eq $null) {
    throw "Could not locate the manual page break paragraph"
}
$rng = $d.Range($pageBreakPara.Range.Start, $pageBreakPara.Range.End)
$rng.InsertXML((New-FlatOpcBody '<w:p><w:r><w:lastRenderedPageBreak/><w:br w:type="page"/></w:r></w:p>'))
